# Update the employee/meme table: lower-case the names (moving the header
# pairing so "Image URL" now sits in B2 next to "bishal" in A2, and the
# first meme image moves up into B1 next to "Name"), matching the
# re-uploaded source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "meme_images/bishal_meme.png"

$ws.Range("A2").Value = "bishal"
$ws.Range("B2").Value = "Image URL"

$ws.Range("A3").Value = "rahul"
$ws.Range("B3").Value = "meme_images/rahul_meme.png"

$ws.Range("A4").Value = "adrija"
$ws.Range("B4").Value = "meme_images/adrija_meme.png"

$ws.Range("A5").Value = "abhishek"
$ws.Range("B5").Value = "meme_images/abhishek_meme.png"

$ws.Range("A6").Value = "prashant"
$ws.Range("B6").Value = "meme_images/prashant_meme.png"

$ws.Range("A7").Value = "sarwesh"
$ws.Range("B7").Value = "meme_images/sarwesh_meme.png"

$ws.Range("A8").Value = "faheem"
$ws.Range("B8").Value = "meme_images/faheem_meme.png"

$ws.Range("A8").Select()
